$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-modified date for every data row
# (rows 2-223). The automatic update bumps this date by one day
# (serial 45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11) for every row.
for ($r = 2; $r -le 223; $r++) {
    $ws.Cells.Item($r, 3).Value = 45180
}
